# Apply the LoginTestData.xlsx changes:
#  - add a "userDetails" worksheet (positioned after "LoginFunc") with the
#    new user-registration sample data, column widths and e-mail hyperlinks
#  - make "userDetails" the active/selected sheet (was "LoginFunc")
#  - drop the stray trailing blank row (+ its customFormat) on "devTestLogin"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. devTestLogin (sheet1): remove the trailing empty row 11 and the
#    customFormat="1" markers that were left on rows 2-10.
# ---------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("devTestLogin")
for ($r = 2; $r -le 10; $r++) {
    $wsLogin.Rows.Item($r).ClearFormats()
}
$wsLogin.Rows.Item(11).Delete()

# ---------------------------------------------------------------------
# 2. Add the new "userDetails" worksheet and move it after "LoginFunc"
#    (re-fetch the sheet by name afterwards -- Move() re-seats the old
#    variable to whatever sheet now sits at its original index).
# ---------------------------------------------------------------------
$wsNew = $wb.Worksheets.Add()
$wsNew.Name = "userDetails"
$wsNew.Move($null, $wb.Worksheets.Item("LoginFunc"))
$ws3 = $wb.Worksheets.Item("userDetails")

# Header row
$ws3.Range("A1").Value = "First  Name "
$ws3.Range("B1").Value = "Last Name"
$ws3.Range("C1").Value = "Email ID"
$ws3.Range("D1").Value = "Phone Number"
$ws3.Range("E1").Value = "User Name"
$ws3.Range("F1").Value = "Password"
$ws3.Range("G1").Value = "Confirm Password"
$ws3.Range("H1").Value = "Login user"
$ws3.Range("I1").Value = "Login password"

# Row 2
$ws3.Range("A2").Value = "Sujith123"
$ws3.Range("B2").Value = "css"
$ws3.Range("C2").Value = "sujith@gmail.com"
$ws3.Range("D2").Value = 123456
$ws3.Range("E2").Value = "jith"
$ws3.Range("F2").Value = "userjith"
$ws3.Range("G2").Value = "userjith"
$ws3.Range("H2").Value = "admin"
$ws3.Range("I2").Value = "useradmin"

# Row 3
$ws3.Range("A3").Value = "Sujith456"
$ws3.Range("B3").Value = "cs123"
$ws3.Range("C3").Value = "sujith@gmail.com"
$ws3.Range("D3").Value = 1234566
$ws3.Range("E3").Value = "jith1"
$ws3.Range("F3").Value = "userjith"
$ws3.Range("G3").Value = "userjith"
$ws3.Range("H3").Value = "admin"
$ws3.Range("I3").Value = "useradmin"

# e-mail hyperlinks on C2/C3 (also paints the built-in Hyperlink style)
$ws3.Hyperlinks.Add($ws3.Range("C2"), "mailto:sujith@gmail.com") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "mailto:sujith@gmail.com") | Out-Null

# Row 4: a lone formatted-but-empty cell under the hyperlink column
$ws3.Range("C4").Value = ""
$ws3.Range("C4").Style = "Hyperlink"

# Column widths (approximate the author's manual "best fit" sizing)
$ws3.Columns.Item(1).ColumnWidth = 9.8
$ws3.Columns.Item(2).ColumnWidth = 8.62
$ws3.Columns.Item(3).ColumnWidth = 14.98
$ws3.Columns.Item(4).ColumnWidth = 12.62
$ws3.Columns.Item(5).ColumnWidth = 11.44
$ws3.Columns.Item(6).ColumnWidth = 7.98
$ws3.Columns.Item(7).ColumnWidth = 15.26
$ws3.Columns.Item(8).ColumnWidth = 8.44
$ws3.Columns.Item(9).ColumnWidth = 12.98

# Make "userDetails" the active/selected sheet & cell (was "LoginFunc")
$ws3.Activate()
$ws3.Range("C6").Select() | Out-Null
